$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Extend the text describing the work done on row 18 (keeps its shared-string slot
# right after the "X" entry).
$ws.Range("E18").Value = "Finished schematic (ready for big review), installed piqt in pycharm, talked with S3 project groups about their proces. Made a knob in pyqt that controls RPI hardware PWM pin."

# Translate / update the Dutch activity descriptions (rows 2-6) to English.
$ws.Range("E2").Value = "Intro conversation with Jeroen Veen"
$ws.Range("E3").Value = "kickoff +  worked on POA"
$ws.Range("E4").Value = "Worked on POA, conversation with jeroen about product functionality, details"
$ws.Range("E5").Value = "Worked on global design, kesselring methode, conversation jeroen about global design."
$ws.Range("E6").Value = "Pros and cons list created. State diagram created. Systeem architectuur enhanced after conversation with jeroen"

# Fill in the missing end time on row 18 (D18), copying the time format from C18.
$ws.Range("D18").Value = 0.625
$ws.Range("C18").Copy()
$ws.Range("D18").PasteSpecial(-4122)

# Add the new log entry on row 19.
$ws.Range("B19").Value = 44824
$ws.Range("B18").Copy()
$ws.Range("B19").PasteSpecial(-4122)

$ws.Range("C19").Value = 0.39583333333333331
$ws.Range("C18").Copy()
$ws.Range("C19").PasteSpecial(-4122)

$ws.Range("E19").Value = "Tidied up component boxes and sorted components for extra efficiency and component ESD safety. Attended workshop (20 min) enh connectors schematic"

$excel.CutCopyMode = 0

# Update the active selection to E19, matching the saved view state.
$ws.Range("E19").Select()
